$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 (Ano 2025) values with the new figures
$ws.Range("B7").Value = 2722510.71
$ws.Range("C7").Value = -38.72467157971106
$ws.Range("D7").Value = 2759
$ws.Range("E7").Value = 2759
$ws.Range("F7").Value = 986.7744508880029
$ws.Range("G7").Value = 5.183021166541657
